$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''41.256.63'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +3.47%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.252.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +2.68%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.03%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''302.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +3.20%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''91.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +5.29%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +2.47%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  +0.06%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  +3.69%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''53.97'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +8.75%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''32.02'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +7.13%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.0793'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +2.35%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  +3.27%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''6.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +3.26%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''2.602.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +2.54%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''14.12'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +3.64%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''2.283.59'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +2.49%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''  +4.10%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''41.203.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +3.63%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '''  +8.51%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = '''  +2.74%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''5.90'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +3.04%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''66.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +2.81%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''239.66'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +1.37%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  +5.04%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -0.20%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  +3.17%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''23.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +5.62%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = '''Cosmos'
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = '''9.64'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +5.69%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = '''Toncoin'
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = '''2.12'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -1.55%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''158.06'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +0.72%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''33.82'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +8.44%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  +0.07%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''5.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +6.49%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.0735'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +4.25%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''3.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +8.65%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''2.36'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +1.54%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''16.56'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +8.97%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '''  +2.81%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '''  +6.46%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  +6.15%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''3.96'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +6.67%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''20.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +17.82%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''2.064.65'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -2.07%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  +3.89%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''10.12'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +5.40%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''2.97'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +12.88%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  -1.34%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''2.474.21'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +2.45%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  +2.90%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  +3.58%  '
$ws.Range("E51").Style = "Normal"

Write-Host "Applied all changes"
